$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F2 276->278, F4 1052->1055, F5 559->562
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 278
$ws1.Range("F4").Value = 1055
$ws1.Range("F5").Value = 562

# Sheet "全部类型" (sheet4): F2 276->278, F4 1052->1055, F6 559->562
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 278
$ws4.Range("F4").Value = 1055
$ws4.Range("F6").Value = 562
